# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ B = 0.04271373187048222; C = 0.306821227259698;  D = 0.1494219747398047; E = 0.4942365360607697 }
    3  = @{ B = 0.6606524410359556;  C = 0.306821227259698;  D = 0.7527432677738641; E = 0.4942365360607697 }
    4  = @{ B = 3.286832544864788;   C = 1.655778082260271;  D = 0.7527432677738641; E = 0.4942365360607697 }
    5  = @{ B = 0.04271373187048222; C = 0.04071648406533734; D = 0.7527432677738641; E = 10.19245300693656 }
    6  = @{ B = 0.6606524410359556;  C = 1.655778082260271;  D = 0.1494219747398047; E = 0.4942365360607697 }
    7  = @{ B = 0.6606524410359556;  C = 10.34677158129881;  D = 0.7527432677738641; E = 10.19245300693656 }
    8  = @{ B = 3.286832544864788;   C = 1.655778082260271;  D = 0.7527432677738641; E = 0.4942365360607697 }
    9  = @{ B = 3.286832544864788;   C = 1.655778082260271;  D = 3.537761648806719;  E = 0.4942365360607697 }
    10 = @{ B = 3.286832544864788;   C = 1.655778082260271;  D = 0.1494219747398047; E = 0.4942365360607697 }
    11 = @{ B = 3.286832544864788;   C = 1.655778082260271;  D = 0.7527432677738641; E = 0.4942365360607697 }
    12 = @{ B = 3.286832544864788;   C = 1.655778082260271;  D = 0.7527432677738641; E = 0.4942365360607697 }
    13 = @{ B = 3.286832544864788;   C = 1.655778082260271;  D = 22.3905356188092;   E = 10.19245300693656 }
}

foreach ($r in $data.Keys) {
    $row = $data[$r]
    $b = $row.B
    $c = $row.C
    $d = $row.D
    $e = $row.E

    $ws.Range("B$r").Value = $b
    $ws.Range("C$r").Value = $c
    $ws.Range("D$r").Value = $d
    $ws.Range("E$r").Value = $e
    $ws.Range("G$r").Value = $b + $c + $d + $e
}
